$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The balance table (iebaltab output) gained two extra variable rows
# ("total_personas" and "ingreso"), each with its own Mean/(SE) row pair,
# inserted right before the trailing footnote row. Insert 4 blank rows at
# row 32 so the footnote (old row 32) is pushed down to row 36.
$ws.Rows("32:35").Insert()

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 32: total_personas - group means / Ns / mean difference
Set-TextValue 32 1 "total_personas"
Set-TextValue 32 2 "880"
Set-TextValue 32 3 "3.526"
Set-TextValue 32 4 "16"
Set-TextValue 32 5 "4.016"
Set-TextValue 32 6 "896"
Set-TextValue 32 7 "0.490"

# Row 33: standard errors for total_personas
Set-TextValue 33 3 "(0.059)"
Set-TextValue 33 5 "(0.097)"

# Row 34: ingreso - group means / Ns / mean difference
Set-TextValue 34 1 "ingreso"
Set-TextValue 34 2 "880"
Set-TextValue 34 3 "2.277"
Set-TextValue 34 4 "16"
Set-TextValue 34 5 "3.102"
Set-TextValue 34 6 "896"
Set-TextValue 34 7 "0.825**"

# Row 35: standard errors for ingreso
Set-TextValue 35 3 "(0.044)"
Set-TextValue 35 5 "(0.103)"

# Row 36 (shifted footnote): update the recorded Stata command to include
# the two new variables that were added to the iebaltab call.
$ws.Cells.Item(36, 1).Value = "If the table includes missing values (.n, .o, .v etc.) see the Missing values section in the help file for the Stata command iebaltab for definitions of these values. Significance: ***=.01, **=.05, *=.1. Full user input as written by user: [iebaltab poblacion_urbana_2009 poblacion_por_localidad_2005 poblacion_2005 personas_por_localidad_2007 personas_por_hogar_2007_localida num_est_transmi icv_2007_localidad gasto_promedio_mensual_2007_loca estrato_mean densidad_urbana_2009 area_urbana_2009 acceso_transmi accesibilidad_arterial accesibilidad_arterial_dummy total_personas ingreso , groupvar(dummy_oxxo) control(0) savexlsx(difmedias_controles_baselines_fixed_2011) replace] "
